$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Sheet is protected; unprotect temporarily so the cell values can be updated
$ws.Unprotect()

# Update the confidential disclaimer text to bump the "as of" date from 2021-03-29 to 2021-03-30
$ws.Range("A7").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-30 for illustrative purposes only and are subject to change."

# Update the model holdings weight/percent-change figures
$ws.Range("D2").Value = 0.8439643094167097
$ws.Range("E2").Value = -0.001310615989515096

$ws.Range("D3").Value = 0.1560356905832903
$ws.Range("E3").Value = 0.003213610586011439

$ws.Range("E4").Value = -0.0006046751714474663

# Restore sheet protection as it was before the edit
$ws.Protect()
